# Update the "想去人数" (interest count) figures on the "展览" and
# "全部类型" worksheets to reflect the latest scrape (gh-pages output
# generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of worksheet name -> { row -> new F value }
$updates = @{
    "展览" = @{
        3  = 522
        4  = 268
        6  = 1121
        11 = 1116
        14 = 760
        19 = 664
        22 = 2148
        25 = 1831
        27 = 2664
        34 = 922
        35 = 1641
        38 = 524
        39 = 133
    }
    "全部类型" = @{
        4  = 522
        5  = 268
        7  = 1121
        12 = 1116
        14 = 760
        23 = 664
        26 = 2148
        31 = 2664
        42 = 922
        43 = 1641
        46 = 524
        47 = 133
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
